# iteration two document update
#
# 1) Drop the last two "Website Additions" bullet items ("Employee
#    check-in, check-out on website" and "Add password encryption for
#    employees in application and website"), keeping "Adding a Reset
#    Password System." as the final bullet in that run.
# 2) Give that remaining bullet's paragraph-mark run explicit sz/szCs
#    (28 half-points = 14pt), matching its sibling bullets.
# 3) Mint the ListLabel82..ListLabel117 character styles (the
#    bullet-level label styles LibreOffice/Word emit alongside the
#    numbering definitions already in numbering.xml) so styles.xml
#    keeps them in step with the document's list levels.

$d = $word.ActiveDocument

# --- 1 & 2: locate the three target bullets by their text and fix them up ---

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Adding a Reset Password System.") {
        $target = $i
        break
    }
}

if ($target -ne $null) {
    # The next two paragraphs are the ones being removed; deleting the
    # paragraph that now sits at $target+1 twice removes both, since
    # each delete shifts the following paragraph up into that slot.
    $after1 = $d.Paragraphs.Item($target + 1)
    $after1.Range.Delete()
    $after2 = $d.Paragraphs.Item($target + 1)
    $after2.Range.Delete()

    # Stamp the paragraph mark of "Adding a Reset Password System." with
    # sz=28 / szCs=28 (14pt), matching the run text formatting already
    # used throughout this bullet list.
    $kept = $d.Paragraphs.Item($target)
    $kept.Range.Font.Size = 14
    $kept.Range.Font.SizeBi = 14
}

# --- 3: mint the ListLabel82..ListLabel117 character styles ---

$s = $d.Styles.Add("ListLabel82", 2)
$s.NameLocal = "ListLabel 82"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s.Font.Size = 14
$s = $d.Styles.Add("ListLabel83", 2)
$s.NameLocal = "ListLabel 83"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel84", 2)
$s.NameLocal = "ListLabel 84"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel85", 2)
$s.NameLocal = "ListLabel 85"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel86", 2)
$s.NameLocal = "ListLabel 86"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel87", 2)
$s.NameLocal = "ListLabel 87"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel88", 2)
$s.NameLocal = "ListLabel 88"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel89", 2)
$s.NameLocal = "ListLabel 89"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel90", 2)
$s.NameLocal = "ListLabel 90"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel91", 2)
$s.NameLocal = "ListLabel 91"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s.Font.Size = 16
$s = $d.Styles.Add("ListLabel92", 2)
$s.NameLocal = "ListLabel 92"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel93", 2)
$s.NameLocal = "ListLabel 93"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel94", 2)
$s.NameLocal = "ListLabel 94"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel95", 2)
$s.NameLocal = "ListLabel 95"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel96", 2)
$s.NameLocal = "ListLabel 96"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel97", 2)
$s.NameLocal = "ListLabel 97"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel98", 2)
$s.NameLocal = "ListLabel 98"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel99", 2)
$s.NameLocal = "ListLabel 99"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel100", 2)
$s.NameLocal = "ListLabel 100"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s.Font.Size = 14
$s = $d.Styles.Add("ListLabel101", 2)
$s.NameLocal = "ListLabel 101"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel102", 2)
$s.NameLocal = "ListLabel 102"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel103", 2)
$s.NameLocal = "ListLabel 103"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel104", 2)
$s.NameLocal = "ListLabel 104"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel105", 2)
$s.NameLocal = "ListLabel 105"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel106", 2)
$s.NameLocal = "ListLabel 106"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel107", 2)
$s.NameLocal = "ListLabel 107"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel108", 2)
$s.NameLocal = "ListLabel 108"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel109", 2)
$s.NameLocal = "ListLabel 109"
$s.QuickStyle = $true
$s.Font.Name = "Calibri"
$s.Font.NameBi = "Symbol"
$s.Font.Size = 14
$s = $d.Styles.Add("ListLabel110", 2)
$s.NameLocal = "ListLabel 110"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel111", 2)
$s.NameLocal = "ListLabel 111"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel112", 2)
$s.NameLocal = "ListLabel 112"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel113", 2)
$s.NameLocal = "ListLabel 113"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel114", 2)
$s.NameLocal = "ListLabel 114"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"
$s = $d.Styles.Add("ListLabel115", 2)
$s.NameLocal = "ListLabel 115"
$s.QuickStyle = $true
$s.Font.NameBi = "Symbol"
$s = $d.Styles.Add("ListLabel116", 2)
$s.NameLocal = "ListLabel 116"
$s.QuickStyle = $true
$s.Font.NameBi = "Courier New"
$s = $d.Styles.Add("ListLabel117", 2)
$s.NameLocal = "ListLabel 117"
$s.QuickStyle = $true
$s.Font.NameBi = "Wingdings"

Write-Output "edit applied"
